$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-machine_master")

# Add a new row (33) of machine master data
$ws.Range("A33").Value = 10032
$ws.Range("B33").Value = "Machine 32"
$ws.Range("C33").Value = "F4-30-B9-D4-CD-6F"
$ws.Range("D33").Value = "FB5962911665"
$ws.Range("E33").Value = "192.168.0.358"
$ws.Range("F33").Value = 1001
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = $true
$ws.Range("I33").Value = "superadmin"
$ws.Range("J33").Value = "now()"

# Update the selected cell to match the final state
$ws.Range("J29").Select()
